# "Doing Updates for Financials" - refresh yearly financial figures on the ALBKY sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(2821100, 2994300, 3335400, 3643100, 3940800, 3706400, 4173000)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $vals[$i]
}

$vals = @(-111900, -109300, -115100, -105800, -103400, -207300, -104400)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(15, 4 + $i).Value = $vals[$i]
}

$vals = @(1769600, 2146000, 4558600, 3561900, 4238600, 4024300, 3439900)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(17, 4 + $i).Value = $vals[$i]
}

$vals = @(1051500, 848400, -1223200, 81200, -297800, -317800, 733100)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(18, 4 + $i).Value = $vals[$i]
}

$vals = @(-866200, -859200, -1069600, -1216000, 2853800, -1190800, -6042000)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(20, 4 + $i).Value = $vals[$i]
}

$vals = @(388800, 161700, -2137700, -968700, 2724900, -1347500, -5204300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(21, 4 + $i).Value = $vals[$i]
}

$vals = @(185300, -10900, -2292700, -1134700, 2556100, -1508600, -5308900)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(23, 4 + $i).Value = $vals[$i]
}

$vals = @(84800, -32800, -905200, -781400, -786700, -288300, -1034200)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(24, 4 + $i).Value = $vals[$i]
}

$vals = @(100400, 21900, -1387500, -353300, 3342800, -1220300, -4274700)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(26, 4 + $i).Value = $vals[$i]
}

$vals = @(100400, 21700, -1387800, -353500, 3342700, -1220700, -4359300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(27, 4 + $i).Value = $vals[$i]
}

$vals = @(-76800, 25500, -151200, -16600, -64100, 6600)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(29, 4 + $i).Value = $vals[$i]
}

$vals = @(866200, 859200, 1069600, 1216000, -2853800, 1190800, 6042000)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(32, 4 + $i).Value = $vals[$i]
}

$vals = @(23600, 47300, -1539000, -370000, 3278600, -1214100, -4359300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(33, 4 + $i).Value = $vals[$i]
}

$vals = @(23600, 47300, -1539000, -370000, 3278600, -1214100, -4359300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(35, 4 + $i).Value = $vals[$i]
}

$vals = @(3704200, 3903400, 6378600, 5354800, 6630000, 3202000, 2360200)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(41, 4 + $i).Value = $vals[$i]
}

$vals = @(717200, 3732300, 5713300, 6137200, 12143800, 2004400, 788900)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(42, 4 + $i).Value = $vals[$i]
}

$vals = @(21200, 24500, 51400, 52000, 112300, 83700, 50300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(47, 4 + $i).Value = $vals[$i]
}

$vals = @(1472500, 1579800, 1665700, 1851900, 3147600, 1354600, 1369900)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(48, 4 + $i).Value = $vals[$i]
}

$vals = @(437400, 416600, 387300, 371900, 545100, 612600, 203700)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(49, 4 + $i).Value = $vals[$i]
}

$vals = @(5183100, 5771800, 5678700, 4200100, 6270400, 2034100, 1645900)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(52, 4 + $i).Value = $vals[$i]
}

$vals = @(68231600, 72786000, 77751100, 81832900, 82687600, 65359800, 66363500)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(54, 4 + $i).Value = $vals[$i]
}

$vals = @(121500, 105400, 111500, 138400, 280500, 126300, 57800)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(59, 4 + $i).Value = $vals[$i]
}

$vals = @(735500, 692100, 449600, 1709400, 878400, 821600, 602800)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(61, 4 + $i).Value = $vals[$i]
}

$vals = @(626400, 487800, 480100, 385500, 881300, 676500, 578700)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(62, 4 + $i).Value = $vals[$i]
}

$vals = @(57462600, 62584400, 67620500, 73212300, 73325600, 64534400, 64773400)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(66, 4 + $i).Value = $vals[$i]
}

$vals = @(2109300, 1054700, 1054700)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(70, 8 + $i).Value = $vals[$i]
}

$vals = @(-1886200, -2439900, -2511000, -1163600, -261500, -3668800, -2738400)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(72, 4 + $i).Value = $vals[$i]
}

$vals = @(10769000, 10201600, 10130700, 8620600, 7252700, -229400, 535400)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(76, 4 + $i).Value = $vals[$i]
}

$vals = @(23600, 47300, -1539000, -370000, 3278600, -1214100, -4359300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(81, 4 + $i).Value = $vals[$i]
}

$vals = @(172200, 154700, 165600, 168500, 160700, 104400, "NA")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(83, 4 + $i).Value = $vals[$i]
}

$vals = @(-2811200, -1466800, -1483100, -4973100, 2923900, -810500, -1268300)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(89, 4 + $i).Value = $vals[$i]
}

$vals = @(-208700, -118400, -89300, -142600, -69700, "NA")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(91, 4 + $i).Value = $vals[$i]
}

$vals = @(2221700, -39300, 578800, 3200400, -1371700, -888800, "NA")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(94, 4 + $i).Value = $vals[$i]
}

$vals = @(-700, -2600, -3200, -108100, -66400)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(96, 6 + $i).Value = $vals[$i]
}

$vals = @(228500, 1660000, 1146600, 500700, -536800, -491900, "NA")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(100, 4 + $i).Value = $vals[$i]
}

$vals = @(-35300, -3700, 5700, -3700, -1100, 8600, "NA")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(101, 4 + $i).Value = $vals[$i]
}

$vals = @(-396300, 150200, 248000, -1275600, 1014300, -2182500, -3406200)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(102, 4 + $i).Value = $vals[$i]
}
